$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Run Date in M1 (45937 -> 45938)
$ws.Range("M1").Value = 45938

# Swap row 271 <-> row 272
$ws.Range("B271").Value = 64973
$ws.Range("B272").Value = 48706
$ws.Range("E271").Value = 35.4
$ws.Range("E272").Value = 39.8
$ws.Range("F271").Value = 150
$ws.Range("F272").Value = -144
$ws.Range("G271").Value = 4995
$ws.Range("G272").Value = -4795.2

# Swap row 309 <-> row 310
$ws.Range("B309").Value = 63565
$ws.Range("B310").Value = 61610
$ws.Range("E309").Value = 109.19
$ws.Range("E310").Value = 122.71
$ws.Range("F309").Value = 60
$ws.Range("F310").Value = -58
$ws.Range("G309").Value = 6162.6
$ws.Range("G310").Value = -5957.18

# Swap row 338 <-> row 339
$ws.Range("B338").Value = 55373
$ws.Range("B339").Value = 63520
$ws.Range("E338").Value = 163.62
$ws.Range("E339").Value = 153.4
$ws.Range("F338").Value = -94
$ws.Range("F339").Value = 97
$ws.Range("G338").Value = -13562.32
$ws.Range("G339").Value = 13995.16

# Swap row 364 <-> row 365
$ws.Range("B364").Value = 57885
$ws.Range("B365").Value = 63652
$ws.Range("E364").Value = 62.28
$ws.Range("E365").Value = 55.42
$ws.Range("F364").Value = 4
$ws.Range("F365").Value = 250
$ws.Range("G364").Value = 208.52
$ws.Range("G365").Value = 13032.5

# Swap row 367 <-> row 368
$ws.Range("B367").Value = 61605
$ws.Range("B368").Value = 63563
$ws.Range("E367").Value = 133.78
$ws.Range("E368").Value = 119.04
$ws.Range("F367").Value = -13
$ws.Range("F368").Value = 15
$ws.Range("G367").Value = -1455.48
$ws.Range("G368").Value = 1679.4

# Swap row 374 <-> row 375
$ws.Range("B374").Value = 60325
$ws.Range("B375").Value = 63560
$ws.Range("E374").Value = 151.57
$ws.Range("E375").Value = 134.87
$ws.Range("F374").Value = -102
$ws.Range("F375").Value = 104
$ws.Range("G374").Value = -12939.72
$ws.Range("G375").Value = 13193.44

# Swap row 381 <-> row 382
$ws.Range("B381").Value = 62865
$ws.Range("B382").Value = 57817
$ws.Range("F381").Value = 151
$ws.Range("F382").Value = 3
$ws.Range("G381").Value = 12051.31
$ws.Range("G382").Value = 239.43

# Swap row 392 <-> row 393
$ws.Range("B392").Value = 57835
$ws.Range("B393").Value = 62933
$ws.Range("F392").Value = 1
$ws.Range("F393").Value = 146
$ws.Range("G392").Value = 59.13
$ws.Range("G393").Value = 8632.98

# Swap row 411 <-> row 412
$ws.Range("B411").Value = 63007
$ws.Range("B412").Value = 57856
$ws.Range("F411").Value = 984
$ws.Range("F412").Value = 2
$ws.Range("G411").Value = 168588.72
$ws.Range("G412").Value = 342.66

# Swap row 423 <-> row 424
$ws.Range("B423").Value = 53082
$ws.Range("B424").Value = 63102
$ws.Range("C423").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("C424").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F423").Value = 1
$ws.Range("F424").Value = 36
$ws.Range("G423").Value = 59.47
$ws.Range("G424").Value = 2140.92

# Swap row 528 <-> row 529
$ws.Range("B528").Value = 47097
$ws.Range("B529").Value = 58047
$ws.Range("D528").Value = 112.28
$ws.Range("D529").Value = 105.54
$ws.Range("E528").Value = 134.16
$ws.Range("E529").Value = 126.1
$ws.Range("F528").Value = 15
$ws.Range("F529").Value = 54
$ws.Range("G528").Value = 1684.2
$ws.Range("G529").Value = 5699.16

# Swap row 571 <-> row 572
$ws.Range("B571").Value = 65069
$ws.Range("B572").Value = 53757
$ws.Range("E571").Value = 14.3
$ws.Range("E572").Value = 16.08
$ws.Range("F571").Value = 172
$ws.Range("F572").Value = -159
$ws.Range("G571").Value = 2313.4
$ws.Range("G572").Value = -2138.55

# Swap row 578 <-> row 579
$ws.Range("B578").Value = 64915
$ws.Range("B579").Value = 45695
$ws.Range("E578").Value = 20.98
$ws.Range("E579").Value = 23.58
$ws.Range("F578").Value = 40
$ws.Range("F579").Value = -36
$ws.Range("G578").Value = 789.2
$ws.Range("G579").Value = -710.28

# Swap row 585 <-> row 586
$ws.Range("B585").Value = 45718
$ws.Range("B586").Value = 64927
$ws.Range("E585").Value = 19.38
$ws.Range("E586").Value = 17.26
$ws.Range("F585").Value = -294
$ws.Range("F586").Value = 295
$ws.Range("G585").Value = -4768.68
$ws.Range("G586").Value = 4784.9

# Swap row 593 <-> row 594
$ws.Range("B593").Value = 64919
$ws.Range("B594").Value = 45702
$ws.Range("E593").Value = 27.97
$ws.Range("E594").Value = 31.43
$ws.Range("F593").Value = 224
$ws.Range("F594").Value = -215
$ws.Range("G593").Value = 5891.2
$ws.Range("G594").Value = -5654.5

# Swap row 679 <-> row 680
$ws.Range("B679").Value = 64810
$ws.Range("B680").Value = 53319
$ws.Range("E679").Value = 291.22
$ws.Range("E680").Value = 310.64
$ws.Range("F679").Value = 7
$ws.Range("F680").Value = -6
$ws.Range("G679").Value = 1917.44
$ws.Range("G680").Value = -1643.52

# Swap row 701 <-> row 702
$ws.Range("B701").Value = 60025
$ws.Range("B702").Value = 64833
$ws.Range("E701").Value = 37.22
$ws.Range("E702").Value = 34.9
$ws.Range("F701").Value = -98
$ws.Range("F702").Value = 99
$ws.Range("G701").Value = -3217.34
$ws.Range("G702").Value = 3250.17

# Swap row 712 <-> row 713
$ws.Range("B712").Value = 60022
$ws.Range("B713").Value = 64830
$ws.Range("E712").Value = 37.22
$ws.Range("E713").Value = 34.9
$ws.Range("F712").Value = -113
$ws.Range("F713").Value = 117
$ws.Range("G712").Value = -3709.79
$ws.Range("G713").Value = 3841.11

# Swap row 864 <-> row 865
$ws.Range("B864").Value = 54751
$ws.Range("B865").Value = 65079
$ws.Range("E864").Value = 46.34
$ws.Range("E865").Value = 43.44
$ws.Range("F864").Value = -19
$ws.Range("F865").Value = 21
$ws.Range("G864").Value = -776.53
$ws.Range("G865").Value = 858.27
